$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Tema" column (B3:B12) values with Tema1..Tema10,
# removing the need for the old "Deportes"/"Política"/"Sociedad" strings.
$temas = @("Tema1","Tema2","Tema3","Tema4","Tema5","Tema6","Tema7","Tema8","Tema9","Tema10")
for ($i = 0; $i -lt $temas.Length; $i++) {
    $row = 3 + $i
    $ws.Range("B$row").Value = $temas[$i]
}

# Update the active selection to reflect the edit (cell C3 only).
$ws.Range("C3").Select()
